# ModelParameters.xlsx edit
#
# Commit summary: `createScenarios()` / `Scenario` gain a new
# `stopIfParameterNotFound` argument. The accompanying example workbook
# (this file) was opened in Excel and a new demo parameter row ("foo" /
# "bar" / 2) was added to the "Global" sheet; a handful of stray,
# content-less styled cells left over from earlier edits were cleared out
# at the same time.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Global")
$ws2 = $wb.Worksheets.Item("Aciclovir")

# --- "Global" sheet -------------------------------------------------------

# Drop the stray, value-less styled cells below the real data (D5, B6, D6,
# D10) - they carried no content, only left-over direct formatting.
$ws1.Range("D5").Clear()
$ws1.Range("B6:D6").Clear()
$ws1.Range("D10").Clear()

# Add the new demo parameter row right after the existing data.
$ws1.Range("A3").Value = "foo"
$ws1.Range("B3").Value = "bar"
$ws1.Range("C3").Value = 2
$ws1.Range("C3").NumberFormat = "0.0000"

$ws1.Range("D3").Select() | Out-Null

# --- "Aciclovir" sheet ------------------------------------------------------

# Remove the (unused) direct row/cell formatting on rows 1-2, keeping only
# the numeric formatting on the "Value" column (C).
$ws2.Rows("1:2").ClearFormats()
$ws2.Range("C1:C2").NumberFormat = "0.0000"
